# [Soumya] Add: Sort And Send Email Template
# Update the Email cell for the first contact and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address in H2 (keeps its existing Hyperlink style)
$ws.Range("H2").Value = "soumyadipta0077@gmail.com"

# Move the active selection to H3 (as recorded after the edit)
$ws.Range("H3").Select()
